# Form the consolidated report: fill in the "Absent" (column H) values
# for the rows that were still blank / unset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
